$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing rows 5-36 down to 6-37
$ws.Rows("5:5").Insert()

# The newly inserted row inherits a plain (borderless) style by default;
# restore the bordered "data row" formatting used throughout the table
# by copying it down from the row directly below.
$ws.Range("A6:F6").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 5 with the new test case data
$ws.Range("A5").Value = "Login Page"
$ws.Range("B5").Value = "TC4"
$ws.Range("C5").Value = "Verify Create_account option in login page"
$ws.Range("D5").Value = "Medium "
$ws.Range("E5").Value = "Sanity"

# Update the selected cell to match the saved view state
$null = $ws.Range("E6").Select()
